$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The watch list gained two newer observation dates. Insert two fresh
# columns right after the broker-name column, pushing the existing
# "UN" marker column and the rating-history column two slots to the
# right (B->D, C->E).
$ws.Columns("B:C").Insert()

# Header row: newest date goes first (leftmost), matching the existing
# newest-to-oldest ordering already used by the sheet.
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Neither new date has any recorded rating activity yet, so every
# broker just gets the same "UN" placeholder the sheet already uses
# for "no rating change" in the other date columns.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Keep the new columns the same (cosmetic) width as the column they
# were cloned from.
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14
$ws.Columns("E").ColumnWidth = 7.14
